# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 3453
    $ws.Range("F3").Value = 28
    $ws.Range("F5").Value = 1837
    $ws.Range("F6").Value = 118

    if ($sheetName -eq "展览") {
        $ws.Range("F7").Value = 343
    } else {
        $ws.Range("F8").Value = 343
    }
}
